$d = $word.ActiveDocument

$pairs = @(
    @("885÷4=221, 1", "254÷4=63, 2"),
    @("669÷3=223, 0", "266÷4=66, 2"),
    @("246÷6=41, 0", "216÷6=36, 0"),
    @("796÷7=113, 5", "953÷9=105, 8"),
    @("719÷4=179, 3", "528÷6=88, 0"),
    @("580÷4=145, 0", "260÷5=52, 0"),
    @("978÷6=163, 0", "671÷3=223, 2"),
    @("291÷9=32, 3", "112÷3=37, 1"),
    @("776÷3=258, 2", "317÷8=39, 5"),
    @("821÷8=102, 5", "221÷6=36, 5"),
    @("697÷6=116, 1", "906÷9=100, 6"),
    @("724÷3=241, 1", "196÷4=49, 0"),
    @("816÷9=90, 6", "153÷3=51, 0"),
    @("671÷2=335, 1", "382÷7=54, 4"),
    @("851÷6=141, 5", "338÷5=67, 3"),
    @("180÷9=20, 0", "908÷2=454, 0"),
    @("585÷9=65, 0", "884÷2=442, 0"),
    @("442÷7=63, 1", "576÷6=96, 0"),
    @("671÷8=83, 7", "449÷2=224, 1"),
    @("134÷8=16, 6", "228÷6=38, 0"),
    @("748÷7=106, 6", "811÷5=162, 1"),
    @("956÷4=239, 0", "722÷6=120, 2"),
    @("544÷7=77, 5", "950÷2=475, 0"),
    @("638÷5=127, 3", "808÷7=115, 3"),
    @("715÷2=357, 1", "146÷8=18, 2")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done"
